$wb = $excel.ActiveWorkbook

# --- Worksheet "Eetu Pihamäki": fill in the time-tracking entry for 12.11.2018 (row 29) ---
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Pvm (date) - 12.11.2018
$ws.Range("A29").Value = 43416
# Aloitusklo (start time) - 17:00
$ws.Range("B29").Value = 0.70833333333333337
# Lopetusklo (end time) - 19:40
$ws.Range("C29").Value = 0.81944444444444453
# Sprint
$ws.Range("E29").Value = 4
# Tehtävä (task description)
$ws.Range("F29").Value = "2 h 20 min yritin Markuksen kanssa saada OpenSSL sertifikaattia toimimaan, niin että Windows 10 Pro työasema löytäisi sen OpenSSL:llä. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2012.11.2018.txt"

# Row grew tall enough to show the wrapped note text in full.
$ws.Rows.Item(29).RowHeight = 90

$wb.Save()
